# Updates the cryptos price-ticker sheet in place: refreshes the "Price" (D)
# and "Volume(1h)" (E) columns for each coin row with the latest scraped
# values (GitHub Actions scheduled refresh).
#
# Several "Price" values are purely numeric-looking strings (e.g. "591.50",
# "7.63") that must stay TEXT (matching the sheet's existing inlineStr/text
# storage and preserving trailing zeros) rather than being auto-coerced to
# numbers by Excel. For those cells we force a Text number format ("@")
# before assigning the value so Excel keeps the literal string instead of
# parsing it into a Number (which would silently drop the trailing zero,
# e.g. "591.50" -> 591.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.416.22"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "3.511.70"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.50"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.72"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +5.99%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("D12").Value = "4.110.09"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "3.510.18"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "64.411.61"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.99"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "396.29"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "3.652.20"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.43%  "
$ws.Range("E27").Value = "  +3.38%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.26"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.49"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("E33").Value = "  +6.81%  "
$ws.Range("D34").Value = "3.540.93"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.37"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.37"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.03"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.90"
$ws.Range("E45").Value = "  -3.60%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.17"
$ws.Range("E47").Value = "  -2.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.80"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").Value = "2.378.02"
$ws.Range("E49").Value = "  -3.77%  "
$ws.Range("E50").Value = "  +0.31%  "
